$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 serotype value was "MgAAV1" (same as row 22); correct it to a new
# distinct label "MgAAV2" to fix the MSA/taxonomy hierarchy.
$ws.Range("C23").Value = "MgAAV2"

# Update the visible selection left behind after the edit.
$ws.Range("C13").Select()
